$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.882.18"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.628.61"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.56"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.43"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.256"
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0613"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.860.10"
$ws.Range("E12").Value = "  -0.79%  "
$ws.Range("D13").Value = "1.621.04"
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.57"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "27.881.48"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.51"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.19"
$ws.Range("E23").Value = "  -5.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.91"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0482"
$ws.Range("E31").Value = "  -0.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "1.400.46"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  +8.02%  "
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.50"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.87"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  -0.76%  "
$ws.Range("D47").Value = "1.770.74"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.17"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.102"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("E51").Value = "  -0.38%  "
